$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert a new contribution-log row at row 5 -----------------------
# Shifts the existing rows 5-23 down to 6-24 (dimension grows to D24).
$ws.Rows("5:5").Insert()

# Carry the bordered cell style from the row above into the freshly
# inserted row (values are not touched - only formats).
$ws.Range("B4:D4").Copy()
$ws.Range("B5:D5").PasteSpecial(-4122)

# New entry: Ryan Conyac edited the proposal video.
$ws.Range("B5").Value = 0
$ws.Range("C5").Value = "Ryan Conyac"
$ws.Range("D5").Value = "Edited video for proposal"

# New entry in the row that used to be the first empty row after the
# log (now row 9 post-insert): Ryan helped Younouss with the key
# controller class.
$ws.Range("B9").Value = 1
$ws.Range("C9").Value = "Ryan Conyac"
$ws.Range("D9").Value = "Helped Younouss integrate key controller class"

# Match the author's last selection (cell D9) when the file was saved.
$ws.Range("D9").Select()
